$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03277891718434
$ws.Range("D2").Value = 1.040946397824372
$ws.Range("E2").Value = 1.041521604215283
$ws.Range("F2").Value = 1.05097486132993
$ws.Range("I2").Value = 1.02359499962809
$ws.Range("J2").Value = 1.037907263792016
$ws.Range("K2").Value = 1.043727305393432
$ws.Range("L2").Value = 1.044300882976483
$ws.Range("M2").Value = 1.053727648330375
$ws.Range("N2").Value = 1.016511895717824

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033887469934977
$ws.Range("D3").Value = 1.041948041778461
$ws.Range("E3").Value = 1.042516183855932
$ws.Range("F3").Value = 1.052073023416918
$ws.Range("I3").Value = 1.023504579208684
$ws.Range("J3").Value = 1.038657446785932
$ws.Range("K3").Value = 1.044538974358904
$ws.Range("L3").Value = 1.045105624527215
$ws.Range("M3").Value = 1.054637628989984
$ws.Range("N3").Value = 1.016769671234598

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034605539528363
$ws.Range("D4").Value = 1.042597193449365
$ws.Range("E4").Value = 1.043160767956478
$ws.Range("F4").Value = 1.052784864318794
$ws.Range("I4").Value = 1.023443335729026
$ws.Range("J4").Value = 1.039143086872412
$ws.Range("K4").Value = 1.045064588212532
$ws.Range("L4").Value = 1.045626757672421
$ws.Range("M4").Value = 1.055227109970175
$ws.Range("N4").Value = 1.016936320762674

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034907598564819
$ws.Range("D5").Value = 1.042870340833675
$ws.Range("E5").Value = 1.043431995946438
$ws.Range("F5").Value = 1.053084422705107
$ws.Range("I5").Value = 1.023416932628374
$ws.Range("J5").Value = 1.039347302746772
$ws.Range("K5").Value = 1.045285654064283
$ws.Range("L5").Value = 1.045845940244149
$ws.Range("M5").Value = 1.055475086085017
$ws.Range("N5").Value = 1.017006344425932

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034958326318581
$ws.Range("D6").Value = 1.042916217778877
$ws.Range("E6").Value = 1.043477550663671
$ws.Range("F6").Value = 1.053134737485721
$ws.Range("I6").Value = 1.023412460900047
$ws.Range("J6").Value = 1.03938159457864
$ws.Range("K6").Value = 1.045322777718984
$ws.Range("L6").Value = 1.045882747711316
$ws.Range("M6").Value = 1.055516731639881
$ws.Range("N6").Value = 1.017018099600504

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034609574938314
$ws.Range("D7").Value = 1.042600842299759
$ws.Range("E7").Value = 1.043164391156677
$ws.Range("F7").Value = 1.052788865852407
$ws.Range("I7").Value = 1.023442985510381
$ws.Range("J7").Value = 1.03914581540783
$ws.Range("K7").Value = 1.045067541721403
$ws.Range("L7").Value = 1.045629686015226
$ws.Range("M7").Value = 1.055230422817505
$ws.Range("N7").Value = 1.016937256563269

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033153400147238
$ws.Range("D8").Value = 1.041284696417669
$ws.Range("E8").Value = 1.041857514646188
$ws.Range("F8").Value = 1.051345729942386
$ws.Range("I8").Value = 1.023565007130014
$ws.Range("J8").Value = 1.038160745625785
$ws.Range("K8").Value = 1.044001527689888
$ws.Range("L8").Value = 1.044572763738084
$ws.Range("M8").Value = 1.054035043151073
$ws.Range("N8").Value = 1.016599042748458

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030593260206128
$ws.Range("D9").Value = 1.03897331779385
$ws.Range("E9").Value = 1.039562498016329
$ws.Range("F9").Value = 1.048812380038185
$ws.Range("I9").Value = 1.023759139049229
$ws.Range("J9").Value = 1.036426634509789
$ws.Range("K9").Value = 1.042126226636269
$ws.Range("L9").Value = 1.042713498095249
$ws.Range("M9").Value = 1.051933720943189
$ws.Range("N9").Value = 1.016001940929356

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028890403018466
$ws.Range("D10").Value = 1.037437692101012
$ws.Range("E10").Value = 1.038037801025489
$ws.Range("F10").Value = 1.047129980899864
$ws.Range("I10").Value = 1.023874603872166
$ws.Range("J10").Value = 1.035271721961216
$ws.Range("K10").Value = 1.040878161823052
$ws.Range("L10").Value = 1.04147613454198
$ws.Range("M10").Value = 1.050536279062445
$ws.Range("N10").Value = 1.015603130729733

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028153968821848
$ws.Range("D11").Value = 1.03677400903684
$ws.Range("E11").Value = 1.037378855507853
$ws.Range("F11").Value = 1.046403028840855
$ws.Range("I11").Value = 1.023921309042974
$ws.Range("J11").Value = 1.034771909872482
$ws.Range("K11").Value = 1.04033824558155
$ws.Range("L11").Value = 1.040940854931738
$ws.Range("M11").Value = 1.049931990277615
$ws.Range("N11").Value = 1.015430268366701

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027880561293763
$ws.Range("D12").Value = 1.036527676121871
$ws.Range("E12").Value = 1.037134283165017
$ws.Range("F12").Value = 1.046133237728193
$ws.Range("I12").Value = 1.023938164268129
$ws.Range("J12").Value = 1.03458629849085
$ws.Range("K12").Value = 1.040137772552842
$ws.Range("L12").Value = 1.04074210461495
$ws.Range("M12").Value = 1.049707652809719
$ws.Range("N12").Value = 1.015366033521727

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027939201947363
$ws.Range("D13").Value = 1.03658050683264
$ws.Range("E13").Value = 1.037186736188707
$ws.Range("F13").Value = 1.046191098369844
$ws.Range("I13").Value = 1.023934571058996
$ws.Range("J13").Value = 1.034626110901284
$ws.Range("K13").Value = 1.040180771260392
$ws.Range("L13").Value = 1.040784733773738
$ws.Range("M13").Value = 1.049755768421108
$ws.Range("N13").Value = 1.015379813290174

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028131366079543
$ws.Range("D14").Value = 1.036753643233756
$ws.Range("E14").Value = 1.037358635216957
$ws.Range("F14").Value = 1.046380723104707
$ws.Range("I14").Value = 1.023922712353276
$ws.Range("J14").Value = 1.034756566339011
$ws.Range("K14").Value = 1.040321672869031
$ws.Range("L14").Value = 1.040924424609596
$ws.Range("M14").Value = 1.049913443976721
$ws.Range("N14").Value = 1.015424959224897

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02824978287455
$ws.Range("D15").Value = 1.036860343274567
$ws.Range("E15").Value = 1.037464572983329
$ws.Range("F15").Value = 1.046497587811238
$ws.Range("I15").Value = 1.023915340507297
$ws.Range("J15").Value = 1.034836949687171
$ws.Range("K15").Value = 1.040408497063663
$ws.Range("L15").Value = 1.04101050286887
$ws.Range("M15").Value = 1.050010609301042
$ws.Range("N15").Value = 1.015452771677571

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028939296955285
$ws.Range("D16").Value = 1.037481764982884
$ws.Range("E16").Value = 1.038081559607479
$ws.Range("F16").Value = 1.047178258732104
$ws.Range("I16").Value = 1.023871434968932
$ws.Range("J16").Value = 1.035304898598095
$ws.Range("K16").Value = 1.040914004920491
$ws.Range("L16").Value = 1.041511669984828
$ws.Range("M16").Value = 1.050576400858557
$ws.Range("N16").Value = 1.015614599368302

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029372055375552
$ws.Range("D17").Value = 1.037871902059073
$ws.Range("E17").Value = 1.038468916127406
$ws.Range("F17").Value = 1.047605637489823
$ws.Range("I17").Value = 1.023843013862435
$ws.Range("J17").Value = 1.03559850366848
$ws.Range("K17").Value = 1.041231231715783
$ws.Range("L17").Value = 1.041826174744545
$ws.Range("M17").Value = 1.050931524683973
$ws.Range("N17").Value = 1.01571606292318

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029624564275758
$ws.Range("D18").Value = 1.038099583291895
$ws.Range("E18").Value = 1.038694976022444
$ws.Range("F18").Value = 1.047855068656766
$ws.Range("I18").Value = 1.023826118601784
$ws.Range("J18").Value = 1.03576978483877
$ws.Range("K18").Value = 1.041416313360778
$ws.Range("L18").Value = 1.042009668951994
$ws.Range("M18").Value = 1.051138740811139
$ws.Range("N18").Value = 1.015775228007061

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029710678264679
$ws.Range("D19").Value = 1.038177237202941
$ws.Range("E19").Value = 1.038772077185316
$ws.Range("F19").Value = 1.047940143398595
$ws.Range("I19").Value = 1.023820303822167
$ws.Range("J19").Value = 1.035828191749015
$ws.Range("K19").Value = 1.041479429660769
$ws.Range("L19").Value = 1.042072244023563
$ws.Range("M19").Value = 1.051209409449441
$ws.Range("N19").Value = 1.015795398900245

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029325615369017
$ws.Range("D20").Value = 1.037830031538081
$ws.Range("E20").Value = 1.03842734388042
$ws.Range("F20").Value = 1.04755976844583
$ws.Range("I20").Value = 1.023846096030192
$ws.Range("J20").Value = 1.03556699990182
$ws.Range("K20").Value = 1.041197191252291
$ws.Range("L20").Value = 1.041792426300611
$ws.Range("M20").Value = 1.050893415146862
$ws.Range("N20").Value = 1.01570517858955

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02807477475583
$ws.Range("D21").Value = 1.036702653672456
$ws.Range("E21").Value = 1.03730801000777
$ws.Range("F21").Value = 1.046324876961286
$ws.Range("I21").Value = 1.02392621804881
$ws.Range("J21").Value = 1.034718149329076
$ws.Range("K21").Value = 1.040280178760443
$ws.Range("L21").Value = 1.040883287030517
$ws.Range("M21").Value = 1.049867009078174
$ws.Range("N21").Value = 1.015411665590211

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027289113602919
$ws.Range("D22").Value = 1.035994917354998
$ws.Range("E22").Value = 1.036605336123456
$ws.Range("F22").Value = 1.045549789715972
$ws.Range("I22").Value = 1.023973741382756
$ws.Range("J22").Value = 1.034184681034617
$ws.Range("K22").Value = 1.03970405585322
$ws.Range("L22").Value = 1.040312117006445
$ws.Range("M22").Value = 1.049222374488053
$ws.Range("N22").Value = 1.015226971610363

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027705532513927
$ws.Range("D23").Value = 1.036369998269674
$ws.Range("E23").Value = 1.036977732875265
$ws.Range("F23").Value = 1.045960551229605
$ws.Range("I23").Value = 1.023948818288725
$ws.Range("J23").Value = 1.034467460094397
$ws.Range("K23").Value = 1.040009427812279
$ws.Range("L23").Value = 1.040614863087235
$ws.Range("M23").Value = 1.049564040342087
$ws.Range("N23").Value = 1.015324895597774

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029346599336146
$ws.Range("D24").Value = 1.037848950650665
$ws.Range("E24").Value = 1.038446128214239
$ws.Range("F24").Value = 1.04758049423388
$ws.Range("I24").Value = 1.023844704313072
$ws.Range("J24").Value = 1.035581235016566
$ws.Range("K24").Value = 1.041212572522178
$ws.Range("L24").Value = 1.04180767561857
$ws.Range("M24").Value = 1.050910634963928
$ws.Range("N24").Value = 1.015710096803594

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031254428840189
$ws.Range("D25").Value = 1.039569933292316
$ws.Range("E25").Value = 1.040154880025948
$ws.Range("F25").Value = 1.049466167349458
$ws.Range("I25").Value = 1.02371141742819
$ws.Range("J25").Value = 1.036874740070763
$ws.Range("K25").Value = 1.042610661690446
$ws.Range("L25").Value = 1.043193785761264
$ws.Range("M25").Value = 1.052476358168658
$ws.Range("N25").Value = 1.016156437977593
